$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.692.13"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "2.204.61"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "260.16"
$r.NumberFormat = "General"
$ws.Range("E5").Value = "  +2.18%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "86.15"
$r.NumberFormat = "General"
$ws.Range("E6").Value = "  +13.26%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("E8").Value = "  -0.02%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.594"
$r.NumberFormat = "General"
$ws.Range("E9").Value = "  +0.64%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "45.14"
$r.NumberFormat = "General"
$ws.Range("E10").Value = "  +7.52%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0916"
$r.NumberFormat = "General"
$ws.Range("E11").Value = "  +0.37%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "7.38"
$r.NumberFormat = "General"
$ws.Range("E12").Value = "  +7.58%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "2.536.49"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "2.237.77"
$ws.Range("E16").Value = "  +0.64%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.782"
$r.NumberFormat = "General"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "43.622.37"
$ws.Range("E19").Value = "  +0.47%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "69.80"
$r.NumberFormat = "General"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("E22").Value = "  +7.89%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "230.83"
$r.NumberFormat = "General"
$ws.Range("E23").Value = "  +0.77%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "8.93"
$r.NumberFormat = "General"
$ws.Range("E24").Value = "  -4.07%  "
$ws.Range("E25").Value = "  +0.02%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "3.54"
$r.NumberFormat = "General"
$ws.Range("E26").Value = "  +5.51%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "10.66"
$r.NumberFormat = "General"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.25"
$r.NumberFormat = "General"
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "2.27"
$r.NumberFormat = "General"
$ws.Range("E29").Value = "  +2.50%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "39.13"
$r.NumberFormat = "General"
$ws.Range("E30").Value = "  +0.83%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "174.48"
$r.NumberFormat = "General"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").Value = "  +2.89%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  +1.50%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.0357"
$r.NumberFormat = "General"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("E38").Value = "  +3.94%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "12.47"
$r.NumberFormat = "General"
$ws.Range("E39").Value = "  +0.23%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "2.86"
$r.NumberFormat = "General"
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("E41").Value = "  -0.53%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "63.35"
$r.NumberFormat = "General"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("E43").Value = "  +4.00%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.198"
$r.NumberFormat = "General"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "8.35"
$r.NumberFormat = "General"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "100.40"
$r.NumberFormat = "General"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("E49").Value = "  +0.41%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.439"
$r.NumberFormat = "General"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("E51").Value = "  +3.16%  "
